$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds plain text in this sheet (e.g. "45.549.22" or
# "317.60"), including values that look numeric. Force it to Text format
# first so Excel does not silently convert these into numbers (which would
# drop significant trailing zeros / reformat the string).
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '45.549.22'
$ws.Range("E2").Value = '  +7.22%  '

# Row 3
$ws.Range("D3").Value = '2.387.51'
$ws.Range("E3").Value = '  +4.48%  '

# Row 4
$ws.Range("E4").Value = '  +0.19%  '

# Row 5
$ws.Range("D5").Value = '114.47'
$ws.Range("E5").Value = '  +10.96%  '

# Row 6
$ws.Range("D6").Value = '317.60'
$ws.Range("E6").Value = '  +2.03%  '

# Row 7
$ws.Range("E7").Value = '  +1.71%  '

# Row 8
$ws.Range("E8").Value = '  -0.07%  '

# Row 9
$ws.Range("E9").Value = '  +4.45%  '

# Row 10
$ws.Range("D10").Value = '42.95'
$ws.Range("E10").Value = '  +11.21%  '

# Row 11
$ws.Range("D11").Value = '0.0938'
$ws.Range("E11").Value = '  +4.60%  '

# Row 12
$ws.Range("E12").Value = '  +6.78%  '

# Row 13
$ws.Range("E13").Value = '  +1.76%  '

# Row 14
$ws.Range("E14").Value = '  +4.52%  '

# Row 15
$ws.Range("D15").Value = '15.90'
$ws.Range("E15").Value = '  +4.13%  '

# Row 16
$ws.Range("D16").Value = '2.749.63'
$ws.Range("E16").Value = '  +4.54%  '

# Row 17
$ws.Range("D17").Value = '2.392.74'
$ws.Range("E17").Value = '  +4.84%  '

# Row 18
$ws.Range("D18").Value = '45.526.11'
$ws.Range("E18").Value = '  +7.27%  '

# Row 19
$ws.Range("D19").Value = '7.57'
$ws.Range("E19").Value = '  +3.65%  '

# Row 20
$ws.Range("E20").Value = '  +3.56%  '

# Row 21
$ws.Range("D21").Value = '13.44'
$ws.Range("E21").Value = '  +0.05%  '

# Row 22
$ws.Range("D22").Value = '74.72'
$ws.Range("E22").Value = '  +2.04%  '

# Row 23
$ws.Range("D23").Value = '3.53'
$ws.Range("E23").Value = '  +3.89%  '

# Row 24
$ws.Range("D24").Value = '267.77'
$ws.Range("E24").Value = '  -0.14%  '

# Row 25
$ws.Range("E25").Value = '  +9.30%  '

# Row 26
$ws.Range("E26").Value = '  +0.11%  '

# Row 27
$ws.Range("D27").Value = '7.72'
$ws.Range("E27").Value = '  +10.24%  '

# Row 28
$ws.Range("D28").Value = '11.32'
$ws.Range("E28").Value = '  +5.41%  '

# Row 29
$ws.Range("D29").Value = '2.35'
$ws.Range("E29").Value = '  +2.43%  '

# Row 30
$ws.Range("D30").Value = '22.90'
$ws.Range("E30").Value = '  +2.65%  '

# Row 31
$ws.Range("D31").Value = '38.87'
$ws.Range("E31").Value = '  +8.82%  '

# Row 32
$ws.Range("D32").Value = '0.0976'
$ws.Range("E32").Value = '  +15.38%  '

# Row 33
$ws.Range("D33").Value = '171.98'
$ws.Range("E33").Value = '  +4.74%  '

# Row 34
$ws.Range("D34").Value = '2.99'
$ws.Range("E34").Value = '  +17.57%  '

# Row 35
$ws.Range("B35").Value = 'Kaspa'
$ws.Range("C35").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D35").Value = '0.121'
$ws.Range("E35").Value = '  +8.32%  '

# Row 36
$ws.Range("B36").Value = 'Stellar'
$ws.Range("C36").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D36").Value = '0.131'
$ws.Range("E36").Value = '  +1.37%  '

# Row 37
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").Value = '4.97'
$ws.Range("E37").Value = '  +10.83%  '

# Row 38
$ws.Range("D38").Value = '4.16'
$ws.Range("E38").Value = '  +15.33%  '

# Row 39
$ws.Range("D39").Value = '3.06'
$ws.Range("E39").Value = '  +11.18%  '

# Row 40
$ws.Range("E40").Value = '  +5.92%  '

# Row 41
$ws.Range("D41").Value = '1.73'
$ws.Range("E41").Value = '  +11.48%  '

# Row 42
$ws.Range("D42").Value = '102.18'
$ws.Range("E42").Value = '  -8.59%  '

# Row 43
$ws.Range("E43").Value = '  +6.47%  '

# Row 44
$ws.Range("D44").Value = '71.56'
$ws.Range("E44").Value = '  +1.54%  '

# Row 45
$ws.Range("D45").Value = '13.24'
$ws.Range("E45").Value = '  +10.01%  '

# Row 46
$ws.Range("E46").Value = '  +0.15%  '

# Row 47
$ws.Range("E47").Value = '  +12.89%  '

# Row 48
$ws.Range("D48").Value = '116.18'
$ws.Range("E48").Value = '  +5.52%  '

# Row 49
$ws.Range("D49").Value = '1.64'
$ws.Range("E49").Value = '  +16.78%  '

# Row 50
$ws.Range("D50").Value = '9.37'
$ws.Range("E50").Value = '  +8.35%  '

# Row 51
$ws.Range("D51").Value = '79.10'
$ws.Range("E51").Value = '  +2.62%  '
